# Excel COM-interop script: build the Walgreens data sheet
#  - rename Sheet1 -> WalgreensDataSheet
#  - populate a small Key/Value table (Item -> Eye Drops)
#  - format header row (bold, yellow fill) and put a thin border around
#    every populated cell (A1:B5)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$wsName = $ws.Name

# --- content -------------------------------------------------------------
# Write "Item" first so the shared-strings table comes out in the same
# order the original workbook used (Item, Key, Value, Eye Drops).
$ws.Range("A2").Value = "Item"
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Value"
$ws.Range("B2").Value = "Eye Drops"

# --- formatting ------------------------------------------------------------
# Build the two required cell formats once on a scratch sheet and paste
# them across - this produces the same two "extra" styles (border-only,
# bold+fill+border) that a normal Format Cells pass would, without piling
# up unused intermediate styles.
$helper = $wb.Worksheets.Add()

$dataTemplate = $helper.Range("A2")
$dataTemplate.Borders.LineStyle = 1

$headerTemplate = $helper.Range("A1")
$headerTemplate.Borders.LineStyle = 1
$headerTemplate.Font.Bold = $true
$headerTemplate.Interior.Color = 65535

$main = $wb.Worksheets.Item($wsName)

$dataRange = $main.Range("A2:B5")
$dataTemplate.Copy()
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headerRange = $main.Range("A1:B1")
$headerTemplate.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$excel.DisplayAlerts = $false
$helper.Delete() | Out-Null
$excel.DisplayAlerts = $true

# --- sheet-level touches ----------------------------------------------------
# Re-resolve the sheet after deleting the scratch sheet - the old
# reference goes stale once the workbook's sheet collection changes.
$main = $wb.Worksheets.Item(1)
$main.Name = "WalgreensDataSheet"
$main.Range("B2").Select() | Out-Null
